# Add v1 .stl files and imported component models
# This edit adds a new BOM row (row 21) for an Adafruit LED Holder, shortens two
# existing Amazon hyperlink URLs (for the M2 heat-set insert and M2 standoff rows),
# and swaps which row references which (now-shortened) URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the two existing M2-related hyperlinks (row 19 / row 20) ---
# Row 19 currently links to the AIEX (threaded insert) URL; after the edit it
# should link to the (shortened) HELIFOUNER (standoffs) URL.
# Row 20 currently links to the HELIFOUNER (standoffs) URL; after the edit it
# should link to the (shortened) AIEX (threaded insert) URL.
# We update the existing Hyperlink objects' Address in place (rather than
# deleting/re-adding) so the cell keeps its original style, and then set the
# cell text to match.

$newHelifounerUrl = "https://www.amazon.com/HELIFOUNER-Spacers-Standoffs-Assortment-Tweezers/dp/B09F8TCLRY/"
$newAiexUrl = "https://www.amazon.com/AIEX-Printing-Embedment-Automotive-M2x3x3-5mm/dp/B0B8GN63S2/"

function Get-HyperlinkAt($sheet, $addr) {
    $found = $null
    foreach ($hl in $sheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $found = $hl
        }
    }
    return $found
}

$hl19 = Get-HyperlinkAt $ws '$G$19'
$hl20 = Get-HyperlinkAt $ws '$G$20'

# Update G19's link/text first so the new (shortened) URL string becomes the
# first newly-introduced shared string.
$hl19.Address = $newHelifounerUrl
$ws.Range("G19").Value2 = $newHelifounerUrl

$hl20.Address = $newAiexUrl
$ws.Range("G20").Value2 = $newAiexUrl

# --- Add the new BOM row 21: Adafruit 5mm Chromed Wide Concave Bevel LED Holder ---

# Copy formatting for the name/mfr/part#/qty/price columns from row 20 (the row
# above), since those are plain direct-formatted (non-named-style) cells.
$ws.Range("B20:F20").Copy()
$ws.Range("B21:F21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ledHolderUrl = "https://www.adafruit.com/product/2178"

# Set the vendor-link URL text before the name text, so new shared strings are
# appended in URL-then-name order (matching how the rest of the sheet reads).
$ws.Range("G21").Value2 = $ledHolderUrl
$ws.Hyperlinks.Add($ws.Range("G21"), $ledHolderUrl) | Out-Null
$ws.Range("G21").HorizontalAlignment = -4131  # xlLeft

$ws.Range("B21").Value2 = "5mm Chromed Wide Concave Bevel LED Holder"
$ws.Range("C21").Value2 = "Adafruit "
$ws.Range("D21").Value2 = 2178
$ws.Range("E21").Value2 = 1
$ws.Range("F21").Value2 = 0.99

# --- Update the sheet view to scroll to / select the new row ---
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("G21").Select()

$wb.Save()
